$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '244.67'
$ws.Range("D3").Value = '21.88'
$ws.Range("D4").Value = '5.391'
$ws.Range("D7").Value = '0.8156'
$ws.Range("D8").Value = '0.9441'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '0.01120'
$ws.Range("E9").Value = '8OneONEBestin24h'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1437'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.07390'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '0.03444'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03056'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09413'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").Value = '4.006'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = '0.001589'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = '0.04815'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("D18").Value = '0.005588'
$ws.Range("D19").Value = '0.004157'
$ws.Range("D20").Value = '0.0009862'
$ws.Range("D21").Value = '3.670'
$ws.Range("D22").Value = '6.427'
$ws.Range("D23").Value = '2.171'
$ws.Range("D26").Value = '0.00007003'
$ws.Range("D41").Value = '0.006402'
$ws.Range("D43").Value = '0.002901'
$ws.Range("D44").Value = '0.006648'
$ws.Range("E44").Value = '43LocalTradersLCT'
$ws.Range("D45").Value = '0.00005247'
$ws.Range("D47").Value = '0.8603'
$ws.Range("D48").Value = '0.002717'
$ws.Range("D49").Value = '0.00002101'
